$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Subscript three character used in PEPE price (0.0₃0743 etc.)
$sub3 = [char]0x2083

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '57.542.75'
$ws.Range("E2").Value = '  -2.34%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.563.05'
$ws.Range("E3").Value = '  -3.68%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '520.33'
$ws.Range("E5").Value = '  -0.69%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.69'
$ws.Range("E6").Value = '  -0.25%  '

# Row 7
$ws.Range("E7").Value = '  -0.12%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.560'
$ws.Range("E8").Value = '  -1.64%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.574.99'
$ws.Range("E9").Value = '  -3.45%  '

# Row 10
$ws.Range("E10").Value = '  -5.31%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0998'
$ws.Range("E11").Value = '  -2.80%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.325'
$ws.Range("E12").Value = '  -3.28%  '

# Row 13
$ws.Range("E13").Value = '  -0.46%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.015.03'
$ws.Range("E14").Value = '  -3.54%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '57.520.68'
$ws.Range("E15").Value = '  -2.32%  '

# Row 16
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.652.41'
$ws.Range("E16").Value = '  -0.26%  '

# Row 17
$ws.Range("B17").Value = 'Avalanche'
$ws.Range("C17").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '20.10'
$ws.Range("E17").Value = '  -4.39%  '

# Row 18
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0000133'
$ws.Range("E18").Value = '  -2.78%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '334.31'
$ws.Range("E19").Value = '  -1.54%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.28'
$ws.Range("E20").Value = '  -2.68%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.13'
$ws.Range("E21").Value = '  -2.61%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.24'
$ws.Range("E22").Value = '  -1.95%  '

# Row 23
$ws.Range("E23").Value = '  -0.08%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.57'
$ws.Range("E24").Value = '  +0.17%  '

# Row 25
$ws.Range("E25").Value = '  -0.09%  '

# Row 26
$ws.Range("E26").Value = '  -0.13%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.400'
$ws.Range("E27").Value = '  -4.93%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.669.70'
$ws.Range("E28").Value = '  -3.82%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.94'
$ws.Range("E29").Value = '  -2.90%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = (@("0.0", "0743") -join $sub3)
$ws.Range("E30").Value = '  -7.54%  '

# Row 31
$ws.Range("E31").Value = '  -0.04%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.25'
$ws.Range("E32").Value = '  -6.59%  '

# Row 33
$ws.Range("E33").Value = '  -0.87%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '18.55'
$ws.Range("E34").Value = '  -1.93%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '148.29'
$ws.Range("E35").Value = '  -1.61%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.01'
$ws.Range("E36").Value = '  -3.34%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.14'
$ws.Range("E37").Value = '  -4.44%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.837'
$ws.Range("E38").Value = '  -8.23%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '35.86'
$ws.Range("E39").Value = '  -2.74%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.832'
$ws.Range("E40").Value = '  -4.27%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.44'
$ws.Range("E41").Value = '  -1.38%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.49'
$ws.Range("E42").Value = '  -2.69%  '

# Row 43
$ws.Range("E43").Value = '  -0.13%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '268.78'
$ws.Range("E44").Value = '  -2.51%  '

# Row 45
$ws.Range("E45").Value = '  -0.07%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0953'
$ws.Range("E46").Value = '  -1.62%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.587'
$ws.Range("E47").Value = '  -4.36%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '18.75'
$ws.Range("E48").Value = '  -5.31%  '

# Row 49
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.72'
$ws.Range("E49").Value = '  +0.06%  '

# Row 50
$ws.Range("B50").Value = 'Hedera'
$ws.Range("C50").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0519'
$ws.Range("E50").Value = '  -2.77%  '

# Row 51
$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.972.19'
$ws.Range("E51").Value = '  -4.33%  '
